$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 2 problems: update the topic/category column (D) for a handful of rows.
$ws.Range("D9").Value  = "comb"
$ws.Range("D10").Value = "comb"
$ws.Range("D11").Value = "numb"
$ws.Range("D13").Value = "numb"
$ws.Range("D14").Value = "comb"
$ws.Range("D16").Value = "geo"

$ws.Range("D49").Value = "comb"
$ws.Range("D51").Value = "geo"
$ws.Range("D52").Value = "geo"
$ws.Range("D53").Value = "comb"
$ws.Range("D55").Value = "alg"
$ws.Range("D56").Value = "alg"

$ws.Range("D89").Value = "numb"
$ws.Range("D90").Value = "comb"
$ws.Range("D92").Value = "comb"
$ws.Range("D93").Value = "geo"
$ws.Range("D94").Value = "alg"
$ws.Range("D95").Value = "alg"
$ws.Range("D96").Value = "geo"

# Reset the view back to the top-left cell / default selection.
$ws.Range("A1").Select() | Out-Null
